$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 describes the "5-way Tactile Switch" order - update quantity, price,
# purchaser and remark to reflect the new order that was placed.
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = 17.44
$ws.Range("E12").Value = "Nathaniel"
$ws.Range("F12").Value = "Six 5-way tactile switches."

# E12 previously had no content/formatting applied (same blank style as E10/E11);
# give it the same look as the rest of row 12 (centered, bordered, size-16 font)
# now that it holds a value.
$ws.Range("E12").Font.Size = $ws.Range("D12").Font.Size
$ws.Range("E12").Borders.LineStyle = $ws.Range("D12").Borders.Item(7).LineStyle
$ws.Range("E12").HorizontalAlignment = $ws.Range("D12").HorizontalAlignment

# Scroll the view back to the top of the sheet and move the active selection
# from F13 to B13.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B13").Select()
